$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.531.97"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.600.98"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'539.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.48%  "
$ws.Range("D6").Value = "'141.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'6.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'0.334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "'0.134"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "3.060.18"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "59.450.00"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "'20.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.610.62"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "'341.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "'10.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'7.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").Value = "0.0₃0744"
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +6.34%  "
$ws.Range("D31").Value = "'5.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'18.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "'149.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "'0.846"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "'0.826"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").Value = "'3.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'272.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "'0.600"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "'10.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.0523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "'18.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.75%  "
$ws.Range("D47").Value = "'0.0223"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "1.939.66"
$ws.Range("D49").Value = "'4.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'111.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("E51").Value = "  +1.68%  "
